$d = $word.ActiveDocument

# --- Paragraph 1: update date and title (two runs separated by a line break) ---
$null = $d.Content.Find.Execute('-04.10.24', $true, $false, $false, $false, $false, $true, 1, $false, '-03.10.24', 2)
$null = $d.Content.Find.Execute('Were RNNs All We Needed?', $true, $false, $false, $false, $false, $true, 1, $false, 'Transformers are Expressive, But Are They Expressive Enough for Regression?', 2)

# --- Paragraphs 2-9: replace full paragraph text directly (avoids autocorrect of quotes) ---
$p = $d.Paragraphs(2)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = 'שוב מאמר על הטרנספורמרים אבל קצת שונה מהמאמר הסטנדרטי על LLMs. המאמר הזה מציג חקירה מעמיקה לגבי expressiveness של הטרנספורמרים, תוך בחינה ספציפית של יכולתם בתור משערכי פונקציות אוניברסליים (כאלו שניתן לקרב איתם כל פונקציה חלקה בדיוק נתון). המחברים מאתגרים טענות קיימות לגבי expressiveness של הטרנספורמרים ומספקים הוכחות תיאורטיות ואמפיריות כאחד שתומכים בהשערתם שהטרנספורמרים מתקשים לקרב (לשערך) באופן מדויק פונקציות חלקות.'

$p = $d.Paragraphs(3)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = 'לפני 4 שנים הוכח שהטרנספורמר(האנקודר) מסוגל לשערך כל פונקציה רציפה אם יש בו מספיק שכבות (בלוקים של טרנספורמר). המשפט הוכח לפני כ 4 שנים והוא מראה שהטרנספורמר בעל שכבות מרובות למעשה יודע לשערך ופונקציה קבועה למקוטעין (piecewise constant) ועם הגודל המינימלי של אינטרוול הקביעות (=רזולוציה) δ הינו קטן מדי אז ניתן לשערך באמצעותו כל פונקציה חלקה בכל דיוק.'

$p = $d.Paragraphs(4)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = 'המאמר המסוקר מתמקד במחקר של הרזולוציה δ הנדרשת לשערוך בדיוק נתון של פונקציה חלקה. התרומה התיאורטית המרכזית של המאמר היא משפט 4.1, אשר קובע חסם עליון על גורם הרזולוציה δ עבור שמכיל מאפיינים שונים של פונקציה מקורבת f.'

$p = $d.Paragraphs(5)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = 'משפט זה משמעותי מכמה סיבות:'

$p = $d.Paragraphs(6)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = 'א) הוא קושר ישירות את גורם הרזולוציה δ לנגזרות של f. קשר זה מבהיר מדוע פונקציות חלקות עם נגזרות המשתנות במהירות מהוות אתגר קשה עבור טרנספורמרים.'

$p = $d.Paragraphs(7)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = 'ב) החסם מראה יחס הפוך בין δ לבין הנגזרות החלקיות של הפונקציה. עבור פונקציות עם נגזרות גדולות, δ חייב להיות קטן כדי לשמור על איכות הקירוב. זה אומר בעצם שאנו צריכים יותר שכבות של טרנספורמרים כדי לקרב בדיוק גבוה את f.'

$p = $d.Paragraphs(8)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = 'ג) המונח האקספוננציאלי 1/(p+md) בחסם מצביע על כך שככל שממד הקלט m או ממד האמבדינג d גדלים, גורם הרזולוציה δ חייב לקטון אקספוננציאלית כדי לשמור על אותה איכות קירוב.'

$p = $d.Paragraphs(9)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = 'ד״א המחברים מספקים הוכחה מפורטת למשפט זה, תחילה למקרה החד-ממדי ולאחר מכן בהכללה לממדים גבוהים יותר..'

# --- Append 7 new paragraphs after paragraph 9 (indices 10..16 of "added") ---
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = 'יתר על כן, המחברים מקשרים את התוצאה התיאורטית הזו להשלכות המעשיות על ארכיטקטורות טרנספורמר. הם מראים שמספר השכבות הנדרש לקירוב הולם גדל כ ((O(m(1/δ)^(dm, מה שהופך ללא ישים מבחינה חישובית עבור δ קטן וממד הקלט בגודל בינוני m. כלומר צריך יותר מדי שכבות הטרנספומרים בשביל זה.'

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = 'המחברים ביצעו ניסויים מקיפים על הטרנספורמר כדי להשלים את ממצאיהם התיאורטיים. הם עשו 2 ניסויים עם הבנצ''מרקים הבאים:'

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = 'א) EXPT-I (רגרסיה): בדיקת יכולתם של טרנספורמרים לקרב ישירות פונקציות חלקות.'

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = 'ב) EXPT-II (״סיווג מקוונטט״): בדיקת יכולתם של טרנספורמרים לקרב פונקציות קבועות למקוטעין.'

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = 'התברר כי הטרנספורמרים מתפקדים באופן גרוע משמעותית ב-EXPT-I בהשוואה ל-EXPT-II, שזה תומך בהשערה שהם מתקשים בקירוב פונקציות חלקות.'

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = 'הגדלת מספר השכבות, ראשי מנגנון ה-attention, או ממדי  אמבדינג אינה משפרת באופן משמעותי את הביצועים על פונקציות חלקות. לעומת הטרנספורמרים מצליחים לקרב באופן הולם פונקציות קבועות למקוטעין עם רזולוציה δ לא קטנה במיוחד.'

$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = 'https://arxiv.org/pdf/2402.15478'

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
